$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.850.03'
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").Value = '3.727.89'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.01%  '
$ws.Range("D7").Value = '3.728.20'
$ws.Range("E7").Value = '  -0.87%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '4.349.37'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = '3.700.13'
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").Value = '68.798.39'
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '496.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.53%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("D34").Value = '3.868.03'
$ws.Range("E34").Value = '  -1.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").Value = '3.661.83'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("E39").Value = '  +1.20%  '
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  -2.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.37'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '40.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.25%  '
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = '2.741.73'
$ws.Range("E51").Value = '  -2.90%  '
